$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price + 1h volume change), and three row
# re-rankings (WrappedEther/Polygon/TRON and TrustWalletToken/FraxShare)
# per the upstream GitHub Actions data refresh.
# All of these source cells are stored as text in the workbook (prices
# such as "1.000" / "26.437.77" and padded percentages like "  -2.31%  "),
# so force a text number format before assigning, otherwise Excel would
# auto-coerce numeric-looking strings into real numbers and mangle the
# formatting (e.g. "1.000" -> 1, "260.60" -> 260.6).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.437.77"
$ws.Range("E2").Value = "  -2.31%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.847.23"
$ws.Range("E3").Value = "  -1.93%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "260.60"
$ws.Range("E5").Value = "  -7.90%  "
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5219"
$ws.Range("E7").Value = "  -1.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3232"
$ws.Range("E8").Value = "  -8.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06757"
$ws.Range("E9").Value = "  -4.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.87"
$ws.Range("E10").Value = "  -7.98%  "
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7702"
$ws.Range("E11").Value = "  -6.41%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07730"
$ws.Range("E12").Value = "  -1.06%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.839.30"
$ws.Range("E13").Value = "  -2.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.85"
$ws.Range("E14").Value = "  -2.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.033"
$ws.Range("E15").Value = "  -3.94%  "
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.14"
$ws.Range("E17").Value = "  -3.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007909"
$ws.Range("E19").Value = "  -3.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.459.58"
$ws.Range("E20").Value = "  -2.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.143.77"
$ws.Range("E21").Value = "  +0.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.537"
$ws.Range("E22").Value = "  -5.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.529"
$ws.Range("E23").Value = "  -6.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.927"
$ws.Range("E24").Value = "  -5.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.356"
$ws.Range("E25").Value = "  -2.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "144.00"
$ws.Range("E26").Value = "  -1.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.659"
$ws.Range("E27").Value = "  -0.98%  "
$ws.Range("E28").Value = "  -3.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.62"
$ws.Range("E29").Value = "  -1.40%  "
$ws.Range("E30").Value = "  -5.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.164"
$ws.Range("E31").Value = "  -5.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08786"
$ws.Range("E32").Value = "  -1.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04793"
$ws.Range("E33").Value = "  -2.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.134"
$ws.Range("E34").Value = "  -4.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.859"
$ws.Range("E35").Value = "  -1.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6910"
$ws.Range("E36").Value = "  -7.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.113"
$ws.Range("E37").Value = "  -5.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01792"
$ws.Range("E38").Value = "  -5.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.210"
$ws.Range("E39").Value = "  -8.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4913"
$ws.Range("E40").Value = "  -8.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "113.09"
$ws.Range("E41").Value = "  -3.56%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.192"
$ws.Range("E42").Value = "  -2.24%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8968"
$ws.Range("E43").Value = "  -9.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.798"
$ws.Range("E45").Value = "  -5.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4204"
$ws.Range("E46").Value = "  -9.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1264"
$ws.Range("E47").Value = "  -8.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.098"
$ws.Range("E48").Value = "  -4.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05885"
$ws.Range("E49").Value = "  -1.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "35.40"
$ws.Range("E50").Value = "  -4.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "59.21"
$ws.Range("E51").Value = "  -5.02%  "
